$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Armatiek")

# Insert a new row before the current row 85 ("migrate"), shifting rows 85-133 down to 86-134.
$ws.Rows.Item(85).Insert()
# Remove any formatting/content copied onto the new row by the Insert operation so that
# only the columns we explicitly populate end up with content/styles.
$ws.Rows.Item(85).Clear()

# Populate the new row with the "metamodelnature" property.
$ws.Range("A85").Value = "metamodelnature"
$ws.Range("F85").Value = "MIM;GROUPING"
$ws.Range("T85").Value = "OPENAPI"

# Re-apply the same cell formatting used by the surrounding rows for the columns involved.
$ws.Range("A86").Copy()
$ws.Range("A85").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F86").Copy()
$ws.Range("F85").PasteSpecial(-4122)

$ws.Range("G86").Copy()
$ws.Range("G85").PasteSpecial(-4122)

$ws.Range("M86").Copy()
$ws.Range("M85").PasteSpecial(-4122)

$ws.Range("S86").Copy()
$ws.Range("S85").PasteSpecial(-4122)

$ws.Range("T84").Copy()
$ws.Range("T85").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection to reflect where editing left off.
$ws.Activate()
$ws.Range("T86").Select()
